# Updated symbol list on Tue Dec 20 09:46:07 UTC 2022 with GitHub Actions
#
# Applies the per-cell value changes produced by the latest crypto price
# refresh: a handful of standalone Price (column D) updates, a block of
# rows (10-18) whose Coin/Link/Price/Volume columns shifted down by one
# ranking slot (a new row -- "One" -- was inserted at rank 9, pushing
# WazirX..CoinExToken down one row and dropping the old rank-17 "One"
# row), and a couple of "Worst in 24h" volume-label tweaks.
#
# Numeric-looking Price values must stay stored as text (inline string in
# the original workbook) rather than being auto-coerced to numbers by
# Excel, so each such cell is briefly marked as Text (NumberFormat "@")
# before the value is written, then restored to the Normal style so no
# stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- standalone Price (column D) corrections -------------------------------
Set-TextValue "D2"  "247.56"
Set-TextValue "D3"  "21.91"
Set-TextValue "D4"  "5.374"
Set-TextValue "D8"  "0.8186"
Set-TextValue "D9"  "0.9323"

# --- rows 10-18: ranking shift (new "One" row inserted at rank 9) ----------
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.0005778"
$ws.Range("E10").Value = "9OneONE"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1441"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07492"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03243"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03078"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09310"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.571"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001608"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04750"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- further standalone Price (column D) corrections ------------------------
Set-TextValue "D19" "0.006387"
Set-TextValue "D20" "0.005061"
Set-TextValue "D21" "0.001034"
Set-TextValue "D22" "0.0001499"
Set-TextValue "D24" "2.165"
Set-TextValue "D25" "0.3309"

# --- Volume(1h) label tweak (row 27: AAXToken gains "Worstin24h" tag) -------
$ws.Range("E27").Value = "26AAXTokenAABWorstin24h"

# --- remaining standalone Price (column D) corrections ----------------------
Set-TextValue "D40" "0.03955"
Set-TextValue "D41" "0.006906"
Set-TextValue "D43" "0.003399"
Set-TextValue "D44" "0.008525"
Set-TextValue "D49" "0.1778"
Set-TextValue "D50" "0.00002099"
